# UCLA.xlsx: rewrite the "Marker" column on the Comp controls sheet so
# each row reads "Marker:Fluorochrome" instead of just the fluorochrome,
# and drop the now-orphaned "CD127 Alexa 647" shared string by reusing it
# as "CD127:Alexa 647" (handled implicitly by simply overwriting the cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comp controls")

$ws.Range("B4").Value  = "LIVE GREEN:FITC"
$ws.Range("B5").Value  = "CD197:PE-A"
$ws.Range("B6").Value  = "CD4:PerCP-Cy5-5-A"
$ws.Range("B7").Value  = "CD45RA:PE-Cy7-A"
$ws.Range("B8").Value  = "CD194:PE-Cy7-A"
$ws.Range("B9").Value  = "CD27:PE-Cy7-A"
$ws.Range("B10").Value = "CD11C:PE-Cy7-A"
$ws.Range("B11").Value = "CD196:PE-Cy7-A"
$ws.Range("B12").Value = "CD38:APC-A"
$ws.Range("B13").Value = "CD127:Alexa 647"
$ws.Range("B14").Value = "CD8:APC-H7-A"
$ws.Range("B15").Value = "CD45RO:APC-H7-A"
$ws.Range("B16").Value = "CD20:APC-H7-A"
$ws.Range("B17").Value = "CD3+19+20:APC-H7-A"
$ws.Range("B18").Value = "CD3:V450"
$ws.Range("B19").Value = "HLA-DR:V500"

# Move the active selection on the sheet (it had been sitting on B12).
$ws.Activate()
$ws.Range("B20").Select()
